# Auto-generated cell updates applied by the scheduled market-data runner.
# For each sheet, updates the currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) with refreshed values. Cells whose new value is empty are
# cleared outright so the underlying XML cell is dropped, matching the source data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 106
$ws.Range("J9").Value = 100.666664
$ws.Range("L9").Value = 100.666664
$ws.Range("N9").Value = -438.666664
$ws.Range("H31").Value = 268.66666
$ws.Range("I31").Value = 268.66666
$ws.Range("K31").Value = 805.9999799999999
$ws.Range("M31").Value = -575.9999799999999
$ws.Range("H111").Value = 879.7143
$ws.Range("I111").Value = 824.25
$ws.Range("K111").Value = 2472.75
$ws.Range("M111").Value = 594.25
$ws.Range("H116").Value = 4997
$ws.Range("I116").Value = 3999.3333
$ws.Range("K116").Value = 3999.3333
$ws.Range("M116").Value = -557.3332999999998
$ws.Range("H118").Value = 1174.25
$ws.Range("I118").Value = 1174.25
$ws.Range("K118").Value = 3522.75
$ws.Range("M118").Value = -1865.75
$ws.Range("H129").Value = 1393.2222
$ws.Range("I129").Value = 590.1667
$ws.Range("J129").Value = 2999.3333
$ws.Range("K129").Value = 1770.5001
$ws.Range("L129").Value = 8997.999899999999
$ws.Range("M129").Value = 3229.4999
$ws.Range("N129").Value = -18997.9999
$ws.Range("H132").Value = 1891.75
$ws.Range("I132").Value = 2047.7142
$ws.Range("K132").Value = 6143.142599999999
$ws.Range("M132").Value = -3613.142599999999
$ws.Range("H135").Value = 932.5
$ws.Range("I135").Value = 932.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8392.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -5857.5
$ws.Range("N135").ClearContents()
$ws.Range("H141").Value = 2606.5
$ws.Range("I141").Value = 2606.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7819.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2639.5
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 599.3
$ws.Range("I2").Value = 491.44446
$ws.Range("J2").Value = 1570
$ws.Range("K2").Value = 491.44446
$ws.Range("L2").Value = 1570
$ws.Range("M2").Value = -378.44446
$ws.Range("N2").Value = -1796
$ws.Range("H32").Value = 2239.0386
$ws.Range("I32").Value = 1928.6
$ws.Range("K32").Value = 1928.6
$ws.Range("M32").Value = -1641.6
$ws.Range("H61").Value = 4366.3335
$ws.Range("I61").Value = 4366.3335
$ws.Range("K61").Value = 4366.3335
$ws.Range("M61").Value = -4154.3335
$ws.Range("H63").Value = 989.8333
$ws.Range("I63").Value = 1057.8
$ws.Range("K63").Value = 1057.8
$ws.Range("M63").Value = -371.8
$ws.Range("H66").Value = 989.8333
$ws.Range("I66").Value = 1057.8
$ws.Range("K66").Value = 5289
$ws.Range("M66").Value = -1857
$ws.Range("H110").Value = 2990.7778
$ws.Range("I110").Value = 727.75
$ws.Range("J110").Value = 4801.2
$ws.Range("K110").Value = 727.75
$ws.Range("L110").Value = 4801.2
$ws.Range("M110").Value = 1317.25
$ws.Range("N110").Value = -8891.200000000001
$ws.Range("H116").Value = 599.3
$ws.Range("I116").Value = 491.44446
$ws.Range("J116").Value = 1570
$ws.Range("K116").Value = 491.44446
$ws.Range("L116").Value = 1570
$ws.Range("M116").Value = 1802.55554
$ws.Range("N116").Value = -6158
$ws.Range("H132").Value = 2941.3333
$ws.Range("I132").Value = 2941.3333
$ws.Range("K132").Value = 8823.999899999999
$ws.Range("M132").Value = -6293.999899999999
$ws.Range("H136").Value = 4366.3335
$ws.Range("I136").Value = 4366.3335
$ws.Range("K136").Value = 13099.0005
$ws.Range("M136").Value = -10549.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 599.3
$ws.Range("I3").Value = 491.44446
$ws.Range("J3").Value = 1570
$ws.Range("K3").Value = 491.44446
$ws.Range("L3").Value = 1570
$ws.Range("M3").Value = -377.44446
$ws.Range("N3").Value = -1798
$ws.Range("H94").Value = 1139.75
$ws.Range("I94").Value = 1074
$ws.Range("K94").Value = 1074
$ws.Range("M94").Value = -623
$ws.Range("H99").Value = 4191.5
$ws.Range("I99").Value = 4564.9
$ws.Range("J99").Value = 2324.5
$ws.Range("K99").Value = 4564.9
$ws.Range("L99").Value = 2324.5
$ws.Range("M99").Value = -3066.9
$ws.Range("N99").Value = -5320.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5739.5713
$ws.Range("H50").Value = 20076.385
$ws.Range("H58").Value = 3337.7646
$ws.Range("I58").Value = 1859.8
$ws.Range("J58").Value = 3953.5833
$ws.Range("K58").Value = 1859.8
$ws.Range("L58").Value = 3953.5833
$ws.Range("M58").Value = -1656.8
$ws.Range("N58").Value = -4359.5833
$ws.Range("H107").Value = 688.7778
$ws.Range("I107").Value = 644.1539
$ws.Range("J107").Value = 804.8
$ws.Range("K107").Value = 644.1539
$ws.Range("L107").Value = 804.8
$ws.Range("M107").Value = 1275.8461
$ws.Range("N107").Value = -4644.8
$ws.Range("H136").Value = 3337.7646
$ws.Range("I136").Value = 1859.8
$ws.Range("J136").Value = 3953.5833
$ws.Range("K136").Value = 5579.4
$ws.Range("L136").Value = 11860.7499
$ws.Range("M136").Value = -3029.4
$ws.Range("N136").Value = -16960.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H9").Value = 187.5
$ws.Range("I9").Value = 400
$ws.Range("J9").Value = 116.666664
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 349.999992
$ws.Range("M9").Value = -976
$ws.Range("N9").Value = -797.999992
$ws.Range("H46").Value = 2008.3334
$ws.Range("I46").Value = 1218.75
$ws.Range("K46").Value = 3656.25
$ws.Range("M46").Value = -3565.25
$ws.Range("H68").Value = 2083.1667
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 2100
$ws.Range("K68").Value = 5997
$ws.Range("L68").Value = 6300
$ws.Range("M68").Value = -5186
$ws.Range("N68").Value = -7922
$ws.Range("H71").Value = 2083.1667
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 2100
$ws.Range("K71").Value = 17991
$ws.Range("L71").Value = 18900
$ws.Range("M71").Value = -13935
$ws.Range("N71").Value = -27012
$ws.Range("H104").Value = 1400
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H127").Value = 3890
$ws.Range("J127").Value = 3890
$ws.Range("L127").Value = 11670
$ws.Range("N127").Value = -21590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 124.434784
$ws.Range("I2").Value = 136.2
$ws.Range("J2").Value = 102.375
$ws.Range("K2").Value = 136.2
$ws.Range("L2").Value = 102.375
$ws.Range("M2").Value = -23.19999999999999
$ws.Range("N2").Value = -328.375
$ws.Range("H59").Value = 10000
$ws.Range("I59").Value = 10000
$ws.Range("K59").Value = 10000
$ws.Range("M59").Value = -9417
$ws.Range("H70").Value = 166667660
$ws.Range("I70").Value = 166667660
$ws.Range("K70").Value = 166667660
$ws.Range("M70").Value = -166667390
$ws.Range("H73").Value = 166667660
$ws.Range("I73").Value = 166667660
$ws.Range("K73").Value = 166667660
$ws.Range("M73").Value = -166666724
$ws.Range("H93").Value = 21000
$ws.Range("J93").Value = 21000
$ws.Range("L93").Value = 21000
$ws.Range("N93").Value = -24744
$ws.Range("H132").Value = 2448.5
$ws.Range("I132").Value = 2265
$ws.Range("K132").Value = 6795
$ws.Range("M132").Value = -4265

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16508.875
$ws.Range("I7").Value = 15504.5
$ws.Range("J7").Value = 18182.834
$ws.Range("K7").Value = 15504.5
$ws.Range("L7").Value = 18182.834
$ws.Range("M7").Value = -15392.5
$ws.Range("N7").Value = -18406.834
$ws.Range("H22").Value = 663
$ws.Range("I22").Value = 480.9091
$ws.Range("K22").Value = 480.9091
$ws.Range("M22").Value = -185.9091
$ws.Range("H27").Value = 663
$ws.Range("I27").Value = 480.9091
$ws.Range("K27").Value = 480.9091
$ws.Range("M27").Value = -373.9091
$ws.Range("H55").Value = 186.53847
$ws.Range("I55").Value = 153.875
$ws.Range("K55").Value = 153.875
$ws.Range("M55").Value = 19.125
$ws.Range("H105").Value = 20666.666
$ws.Range("J105").Value = 20666.666
$ws.Range("L105").Value = 20666.666
$ws.Range("N105").Value = -27654.666
$ws.Range("H126").Value = 16508.875
$ws.Range("I126").Value = 15504.5
$ws.Range("J126").Value = 18182.834
$ws.Range("K126").Value = 46513.5
$ws.Range("L126").Value = 54548.50199999999
$ws.Range("M126").Value = -44043.5
$ws.Range("N126").Value = -59488.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10565.444
$ws.Range("I26").Value = 12
$ws.Range("J26").Value = 11884.625
$ws.Range("K26").Value = 12
$ws.Range("L26").Value = 11884.625
$ws.Range("M26").Value = 281
$ws.Range("N26").Value = -12470.625
